# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.297.22'
$ws.Cells.Item(2, 5).Value = '  +1.23%  '
$ws.Cells.Item(3, 4).Value = '1.695.45'
$ws.Cells.Item(3, 5).Value = '  +1.74%  '
$ws.Cells.Item(4, 4).Value = "'" + '1.005'
$ws.Cells.Item(4, 5).Value = '  -0.51%  '
$ws.Cells.Item(5, 4).Value = "'" + '219.89'
$ws.Cells.Item(5, 5).Value = '  +1.00%  '
$ws.Cells.Item(6, 4).Value = "'" + '0.5259'
$ws.Cells.Item(6, 5).Value = '  +4.47%  '
$ws.Cells.Item(7, 4).Value = "'" + '1.005'
$ws.Cells.Item(7, 5).Value = '  -0.50%  '
$ws.Cells.Item(8, 4).Value = "'" + '0.2702'
$ws.Cells.Item(8, 5).Value = '  +3.07%  '
$ws.Cells.Item(9, 4).Value = "'" + '0.06473'
$ws.Cells.Item(9, 5).Value = '  +2.68%  '
$ws.Cells.Item(10, 4).Value = "'" + '22.20'
$ws.Cells.Item(10, 5).Value = '  +4.47%  '
$ws.Cells.Item(11, 4).Value = "'" + '0.07484'
$ws.Cells.Item(11, 5).Value = '  +1.67%  '
$ws.Cells.Item(12, 4).Value = '1.715.56'
$ws.Cells.Item(12, 5).Value = '  +2.64%  '
$ws.Cells.Item(13, 4).Value = "'" + '4.567'
$ws.Cells.Item(13, 5).Value = '  +1.11%  '
$ws.Cells.Item(14, 4).Value = "'" + '0.5888'
$ws.Cells.Item(14, 5).Value = '  +2.93%  '
$ws.Cells.Item(15, 4).Value = "'" + '0.000008622'
$ws.Cells.Item(15, 5).Value = '  +2.71%  '
$ws.Cells.Item(16, 4).Value = "'" + '64.98'
$ws.Cells.Item(16, 5).Value = '  +0.98%  '
$ws.Cells.Item(17, 4).Value = '26.396.65'
$ws.Cells.Item(17, 5).Value = '  +0.95%  '
$ws.Cells.Item(18, 4).Value = "'" + '5.000'
$ws.Cells.Item(18, 5).Value = '  +1.53%  '
$ws.Cells.Item(19, 5).Value = '  -0.30%  '
$ws.Cells.Item(20, 4).Value = "'" + '10.87'
$ws.Cells.Item(20, 5).Value = '  +1.11%  '
$ws.Cells.Item(21, 4).Value = "'" + '191.78'
$ws.Cells.Item(21, 5).Value = '  +2.97%  '
$ws.Cells.Item(22, 4).Value = "'" + '6.270'
$ws.Cells.Item(22, 5).Value = '  +1.88%  '
$ws.Cells.Item(23, 5).Value = '  -0.38%  '
$ws.Cells.Item(24, 4).Value = "'" + '145.47'
$ws.Cells.Item(24, 5).Value = '  +2.03%  '
$ws.Cells.Item(25, 4).Value = "'" + '7.714'
$ws.Cells.Item(25, 5).Value = '  +1.65%  '
$ws.Cells.Item(26, 4).Value = "'" + '0.1241'
$ws.Cells.Item(26, 5).Value = '  +6.62%  '
$ws.Cells.Item(27, 4).Value = "'" + '15.94'
$ws.Cells.Item(27, 5).Value = '  +1.93%  '
$ws.Cells.Item(28, 4).Value = "'" + '0.06839'
$ws.Cells.Item(28, 5).Value = '  +18.84%  '
$ws.Cells.Item(29, 4).Value = "'" + '1.345'
$ws.Cells.Item(29, 5).Value = '  +3.47%  '
$ws.Cells.Item(30, 4).Value = "'" + '1.331'
$ws.Cells.Item(30, 5).Value = '  +0.63%  '
$ws.Cells.Item(31, 4).Value = "'" + '3.623'
$ws.Cells.Item(31, 5).Value = '  +4.10%  '
$ws.Cells.Item(32, 4).Value = "'" + '3.575'
$ws.Cells.Item(32, 5).Value = '  +2.81%  '
$ws.Cells.Item(33, 4).Value = "'" + '1.675'
$ws.Cells.Item(33, 5).Value = '  +1.05%  '
$ws.Cells.Item(34, 4).Value = "'" + '1.034'
$ws.Cells.Item(34, 5).Value = '  +3.44%  '
$ws.Cells.Item(35, 4).Value = "'" + '0.6250'
$ws.Cells.Item(35, 5).Value = '  +5.00%  '
$ws.Cells.Item(36, 4).Value = "'" + '2.385'
$ws.Cells.Item(36, 5).Value = '  +0.40%  '
$ws.Cells.Item(37, 4).Value = "'" + '2.716'
$ws.Cells.Item(37, 5).Value = '  +2.64%  '
$ws.Cells.Item(38, 4).Value = "'" + '6.318'
$ws.Cells.Item(38, 5).Value = '  +6.96%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).Value = "'" + '0.01624'
$ws.Cells.Item(39, 5).Value = '  +1.77%  '
$ws.Cells.Item(40, 2).Value = 'Maker'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(40, 4).Value = '1.106.54'
$ws.Cells.Item(40, 5).Value = '  +2.51%  '
$ws.Cells.Item(41, 4).Value = "'" + '0.8786'
$ws.Cells.Item(41, 5).Value = '  +2.67%  '
$ws.Cells.Item(42, 4).Value = "'" + '1.017'
$ws.Cells.Item(42, 5).Value = '  +0.95%  '
$ws.Cells.Item(43, 4).Value = "'" + '101.05'
$ws.Cells.Item(43, 5).Value = '  +1.72%  '
$ws.Cells.Item(44, 4).Value = '1.844.65'
$ws.Cells.Item(44, 5).Value = '  +1.27%  '
$ws.Cells.Item(45, 4).Value = "'" + '0.00000000115'
$ws.Cells.Item(45, 5).Value = '  +4.20%  '
$ws.Cells.Item(46, 4).Value = "'" + '57.20'
$ws.Cells.Item(46, 5).Value = '  +2.79%  '
$ws.Cells.Item(47, 4).Value = "'" + '8.192'
$ws.Cells.Item(47, 5).Value = '  +1.89%  '
$ws.Cells.Item(48, 5).Value = '  +0.23%  '
$ws.Cells.Item(49, 4).Value = "'" + '0.05266'
$ws.Cells.Item(49, 5).Value = '  +1.47%  '
$ws.Cells.Item(50, 4).Value = "'" + '0.4295'
$ws.Cells.Item(50, 5).Value = '  -0.58%  '
$ws.Cells.Item(51, 4).Value = "'" + '6.052'
$ws.Cells.Item(51, 5).Value = '  +4.41%  '
